# Update wealth tax support statistics with final computed data
# (re-run of prepare & render pipeline produced slightly different
# bootstrap means / CIs for several country-scenario rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 73.8610351305732
$ws.Range("C2").Value = 72.5008958144439
$ws.Range("D2").Value = 75.2211744467024
$ws.Range("C3").Value = 75.6468093962197
$ws.Range("D3").Value = 79.8456880508458
$ws.Range("B11").Value = 72.4772494650826
$ws.Range("C11").Value = 69.1071100070964
$ws.Range("D11").Value = 75.8473889230687
$ws.Range("B12").Value = 77.191348422519
$ws.Range("C12").Value = 72.7701077063205
$ws.Range("D12").Value = 81.6125891387176
$ws.Range("B14").Value = 67.3083606323479
$ws.Range("C14").Value = 64.3592136887325
$ws.Range("D14").Value = 70.2575075759634
$ws.Range("B15").Value = 69.2201376910046
$ws.Range("C15").Value = 67.7922840724785
$ws.Range("D15").Value = 70.6479913095308
$ws.Range("C16").Value = 69.1658439404388
$ws.Range("D16").Value = 73.7340076205935
$ws.Range("B24").Value = 66.1253305020267
$ws.Range("C24").Value = 62.5272644956226
$ws.Range("D24").Value = 69.7233965084307
$ws.Range("B25").Value = 69.496882977473
$ws.Range("C25").Value = 64.6384824294078
$ws.Range("D25").Value = 74.3552835255383
$ws.Range("B27").Value = 66.5790404538871
$ws.Range("C27").Value = 63.6524209290054
$ws.Range("D27").Value = 69.5056599787687
$ws.Range("B28").Value = 68.2230493525959
$ws.Range("C28").Value = 66.7801873237022
$ws.Range("D28").Value = 69.6659113814897
$ws.Range("C29").Value = 69.2399647889905
$ws.Range("D29").Value = 73.8430897241188
$ws.Range("B37").Value = 60.6507362007239
$ws.Range("C37").Value = 56.9227438052566
$ws.Range("D37").Value = 64.3787285961912
$ws.Range("B38").Value = 73.8494994443583
$ws.Range("C38").Value = 68.9888192472245
$ws.Range("D38").Value = 78.7101796414921
$ws.Range("B40").Value = 63.6399135871967
$ws.Range("C40").Value = 60.7151665904094
$ws.Range("D40").Value = 66.5646605839841
